$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bulk-updated
# from 45192 (2023-09-23) to 45202 (2023-10-03) for every data row (2-498).
$oldValue = 45192
$newValue = 45202

$range = $ws.Range("C2:C498")
foreach ($cell in $range.Cells) {
    if ($cell.Value2() -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
